$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 383.14285
$ws.Range("I33").Value = 236.4
$ws.Range("K33").Value = 236.4
$ws.Range("M33").Value = -7.400000000000006
$ws.Range("H64").Value = 6840.7334
$ws.Range("I64").Value = 6225.9165
$ws.Range("K64").Value = 6225.9165
$ws.Range("M64").Value = -5977.9165
$ws.Range("H67").Value = 6840.7334
$ws.Range("I67").Value = 6225.9165
$ws.Range("K67").Value = 6225.9165
$ws.Range("M67").Value = -5367.9165
$ws.Range("H96").Value = 4762619.5
$ws.Range("I96").Value = 7936832
$ws.Range("J96").Value = 1301.5
$ws.Range("K96").Value = 23810496
$ws.Range("L96").Value = 3904.5
$ws.Range("M96").Value = -23809123
$ws.Range("N96").Value = -6650.5
$ws.Range("H98").Value = 1599.4193
$ws.Range("I98").Value = 1638
$ws.Range("J98").Value = 1239.3334
$ws.Range("K98").Value = 1638
$ws.Range("L98").Value = 1239.3334
$ws.Range("M98").Value = -140
$ws.Range("N98").Value = -4235.3334
$ws.Range("H100").Value = 1390.9474
$ws.Range("I100").Value = 1081.0667
$ws.Range("K100").Value = 1081.0667
$ws.Range("M100").Value = -540.0667000000001
$ws.Range("H122").Value = 1599.4193
$ws.Range("I122").Value = 1638
$ws.Range("J122").Value = 1239.3334
$ws.Range("K122").Value = 4914
$ws.Range("L122").Value = 3718.0002
$ws.Range("M122").Value = -2464
$ws.Range("N122").Value = -8618.0002
$ws.Range("H132").Value = 24653.58
$ws.Range("I132").Value = 26009.45
$ws.Range("K132").Value = 78028.35000000001
$ws.Range("M132").Value = -75498.35000000001
$ws.Range("H135").Value = 3674.375
$ws.Range("I135").Value = 3232.5
$ws.Range("K135").Value = 29092.5
$ws.Range("M135").Value = -26557.5
$ws.Range("H138").Value = 2412.2856
$ws.Range("I138").Value = 1692.862
$ws.Range("J138").Value = 3185
$ws.Range("K138").Value = 5078.586
$ws.Range("L138").Value = 9555
$ws.Range("M138").Value = 61.41399999999976
$ws.Range("N138").Value = -19835
$ws.Range("H141").Value = 1831.6666
$ws.Range("I141").Value = 1247.5
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 3742.5
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 1437.5
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6128.9
$ws.Range("I61").Value = 774.35297
$ws.Range("K61").Value = 774.35297
$ws.Range("M61").Value = -562.35297
$ws.Range("H97").Value = 1897.0435
$ws.Range("I97").Value = 1686.421
$ws.Range("K97").Value = 1686.421
$ws.Range("M97").Value = -1190.421
$ws.Range("H102").Value = 4161.154
$ws.Range("I102").Value = 4419.1
$ws.Range("K102").Value = 4419.1
$ws.Range("M102").Value = -2797.1
$ws.Range("H110").Value = 556.1905
$ws.Range("I110").Value = 562.2105
$ws.Range("K110").Value = 562.2105
$ws.Range("M110").Value = 1482.7895
$ws.Range("H136").Value = 6128.9
$ws.Range("I136").Value = 774.35297
$ws.Range("K136").Value = 2323.05891
$ws.Range("M136").Value = 226.9410899999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1047.6061
$ws.Range("I94").Value = 914.6667
$ws.Range("J94").Value = 1280.25
$ws.Range("K94").Value = 914.6667
$ws.Range("L94").Value = 1280.25
$ws.Range("M94").Value = -463.6667
$ws.Range("N94").Value = -2182.25
$ws.Range("H105").Value = 3153.7856
$ws.Range("I105").Value = 2172.6667
$ws.Range("K105").Value = 2172.6667
$ws.Range("M105").Value = -425.6667000000002
$ws.Range("H134").Value = 7726.6665
$ws.Range("I134").Value = 7987.189
$ws.Range("K134").Value = 23961.567
$ws.Range("M134").Value = -21426.567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1960.3846
$ws.Range("J58").Value = 2721.889
$ws.Range("L58").Value = 2721.889
$ws.Range("N58").Value = -3127.889
$ws.Range("H93").Value = 17498.834
$ws.Range("I93").Value = 17498.834
$ws.Range("K93").Value = 17498.834
$ws.Range("M93").Value = -15626.834
$ws.Range("H94").Value = 1973
$ws.Range("I94").Value = 1999
$ws.Range("J94").Value = 1967.8
$ws.Range("K94").Value = 1999
$ws.Range("L94").Value = 1967.8
$ws.Range("M94").Value = -1548
$ws.Range("N94").Value = -2869.8
$ws.Range("H99").Value = 14170.429
$ws.Range("I99").Value = 11396.667
$ws.Range("J99").Value = 16250.75
$ws.Range("K99").Value = 11396.667
$ws.Range("L99").Value = 16250.75
$ws.Range("M99").Value = -9898.666999999999
$ws.Range("N99").Value = -19246.75
$ws.Range("H126").Value = 14170.429
$ws.Range("I126").Value = 11396.667
$ws.Range("J126").Value = 16250.75
$ws.Range("K126").Value = 34190.001
$ws.Range("L126").Value = 48752.25
$ws.Range("M126").Value = -31720.001
$ws.Range("N126").Value = -53692.25
$ws.Range("H132").Value = 56869.168
$ws.Range("I132").Value = 72567.86
$ws.Range("K132").Value = 217703.58
$ws.Range("M132").Value = -215173.58
$ws.Range("H136").Value = 1960.3846
$ws.Range("J136").Value = 2721.889
$ws.Range("L136").Value = 8165.667
$ws.Range("N136").Value = -13265.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1144.6923
$ws.Range("I5").Value = 898
$ws.Range("K5").Value = 2694
$ws.Range("M5").Value = -2582
$ws.Range("H37").Value = 41868.21
$ws.Range("J37").Value = 41868.21
$ws.Range("L37").Value = 125604.63
$ws.Range("N37").Value = -125828.63
$ws.Range("H102").Value = 4340.3335
$ws.Range("I102").Value = 4340.3335
$ws.Range("K102").Value = 13021.0005
$ws.Range("M102").Value = -10587.0005
$ws.Range("H132").Value = 1099.5
$ws.Range("J132").Value = 1199.3334
$ws.Range("L132").Value = 10794.0006
$ws.Range("N132").Value = -15854.0006
$ws.Range("H135").Value = 1144.6923
$ws.Range("I135").Value = 898
$ws.Range("K135").Value = 8082
$ws.Range("M135").Value = -5547

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 579.2
$ws.Range("I2").Value = 797.6667
$ws.Range("J2").Value = 251.5
$ws.Range("K2").Value = 797.6667
$ws.Range("L2").Value = 251.5
$ws.Range("M2").Value = -684.6667
$ws.Range("N2").Value = -477.5
$ws.Range("H24").Value = 7749.5
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H107").Value = 444.44446
$ws.Range("I107").Value = 107.4
$ws.Range("J107").Value = 865.75
$ws.Range("K107").Value = 107.4
$ws.Range("L107").Value = 865.75
$ws.Range("M107").Value = 1812.6
$ws.Range("N107").Value = -4705.75
$ws.Range("H124").Value = 48666.332
$ws.Range("J124").Value = 48666.332
$ws.Range("L124").Value = 48666.332
$ws.Range("N124").Value = -58486.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1210.8695
$ws.Range("I93").Value = 814.7778
$ws.Range("J93").Value = 2636.8
$ws.Range("K93").Value = 814.7778
$ws.Range("L93").Value = 2636.8
$ws.Range("M93").Value = 433.2222
$ws.Range("N93").Value = -5132.8
$ws.Range("H100").Value = 3054.4443
$ws.Range("I100").Value = 2811.25
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 2811.25
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -2270.25
$ws.Range("N100").Value = -6082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1118.375
$ws.Range("I107").Value = 890.8
$ws.Range("J107").Value = 1497.6666
$ws.Range("K107").Value = 2672.4
$ws.Range("L107").Value = 4492.9998
$ws.Range("M107").Value = -752.3999999999996
$ws.Range("N107").Value = -8332.9998
$ws.Range("H132").Value = 1622.1111
$ws.Range("I132").Value = 1099.8334
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 3299.5002
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -769.5001999999999
$ws.Range("N132").Value = -13060.0001
